$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SQL-style escaping of apostrophes (' -> '') in the flavor_text / effect
# columns for the contest-effect rows whose text contains an apostrophe.
$ws.Cells.Item(11, 4).Value = "Startles the Pokémon that has the judge''s attention."
$ws.Cells.Item(13, 4).Value = "Works well if it''s the same type as the one before."
$ws.Cells.Item(13, 5).Value = "If the last Pokémon''s appeal is the same type as this move, user earns six points instead of two."
$ws.Cells.Item(22, 5).Value = "Shuffles the next turn''s turn order."
$ws.Cells.Item(23, 4).Value = "Shifts the judge''s attention from others."
$ws.Cells.Item(26, 4).Value = "The appeal''s quality depends on its timing."
$ws.Cells.Item(30, 4).Value = "The appeal works well if the user''s condition is good."
$ws.Cells.Item(33, 4).Value = "Ups the user''s condition.  Helps prevent nervousness."
